$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - rows 4..10 column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 252
$ws1.Range("F5").Value = 3099
$ws1.Range("F6").Value = 62
$ws1.Range("F7").Value = 3876
$ws1.Range("F8").Value = 475
$ws1.Range("F9").Value = 975
$ws1.Range("F10").Value = 34

# Sheet "全部类型" (sheet4.xml) - rows 5..11 column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 252
$ws4.Range("F6").Value = 3099
$ws4.Range("F7").Value = 62
$ws4.Range("F8").Value = 3876
$ws4.Range("F9").Value = 475
$ws4.Range("F10").Value = 975
$ws4.Range("F11").Value = 34

$wb.Save()
